$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet by duplicating "2022-Q2" (so header
#    row, column-A styling and sheet formatting match the existing quarter
#    sheets exactly), place it right before "2022-Q2", then rename it and
#    overwrite its contents with the 2022-Q3 figures.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The duplicated sheet has 4 data rows (rows 2-5); the new quarter needs 6
# data rows (rows 2-7), so insert two extra rows copying existing formatting.
$q3.Rows.Item(6).Insert()
$q3.Range("A5:H5").Copy($q3.Range("A6:H6"))
$q3.Rows.Item(7).Insert()
$q3.Range("A6:H6").Copy($q3.Range("A7:H7"))

# Columns B (fund code) and D:G (scale/position/value figures) are stored as
# text in this workbook (e.g. leading zeros in fund codes) - force text so
# Excel doesn't reinterpret the numeric-looking strings as numbers.
$q3.Range("B2:B7").NumberFormat = "@"
$q3.Range("D2:G7").NumberFormat = "@"

# Row 2
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "000369"
$q3.Range("C2").Value = "广发全球医疗保健（QDII）人民币A"
$q3.Range("D2").Value = "2.76"
$q3.Range("E2").Value = "83.19"
$q3.Range("F2").Value = "2.96"
$q3.Range("G2").Value = "0.0817"
$q3.Range("H2").Value = 7

# Row 3
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "000370"
$q3.Range("C3").Value = "广发全球医疗保健（QDII）美元A"
$q3.Range("D3").Value = "2.75"
$q3.Range("E3").Value = "83.19"
$q3.Range("F3").Value = "2.96"
$q3.Range("G3").Value = "0.0814"
$q3.Range("H3").Value = 7

# Row 4
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "006555"
$q3.Range("C4").Value = "浦银安盛全球智能科技股票（QDII）A"
$q3.Range("D4").Value = "0.25"
$q3.Range("E4").Value = "84.65"
$q3.Range("F4").Value = "2.74"
$q3.Range("G4").Value = "0.0068"
$q3.Range("H4").Value = 8

# Row 5
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "016280"
$q3.Range("C5").Value = "广发全球医疗保健（QDII）人民币C"
$q3.Range("D5").Value = "0.02"
$q3.Range("E5").Value = "83.19"
$q3.Range("F5").Value = "2.96"
$q3.Range("G5").Value = "0.0006"
$q3.Range("H5").Value = 7

# Row 6
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "016281"
$q3.Range("C6").Value = "广发全球医疗保健（QDII）美元C"
$q3.Range("D6").Value = "0.02"
$q3.Range("E6").Value = "83.19"
$q3.Range("F6").Value = "2.96"
$q3.Range("G6").Value = "0.0006"
$q3.Range("H6").Value = 7

# Row 7
$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "014002"
$q3.Range("C7").Value = "浦银安盛全球智能科技股票（QDII）C"
$q3.Range("D7").Value = "0.01"
$q3.Range("E7").Value = "84.65"
$q3.Range("F7").Value = "2.74"
$q3.Range("G7").Value = "0.0003"
$q3.Range("H7").Value = 8

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q3 (all
#    other quarters shift down by one row), copying formatting from the row
#    that is about to become row 3 so the styles (bold index column, etc.)
#    line up exactly like the rest of the table.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A3:D3").Copy($total.Range("A2:D2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.17

# Column A is a plain 0-based row counter, independent of the date labels -
# re-number it sequentially for the rows that shifted down (rows 3-7).
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# The quarter-label / count / value columns for the shifted rows already
# carry the right data (each row kept the data of the row above it), so no
# further changes are required there.
